$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Proximity": append rows 53-63 (new door ENTER/EXIT log entries)
# ---------------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

# Column A holds ISO-style date text ("2026-02-01"). Pre-format as Text so
# Excel keeps the literal string instead of auto-converting it to a date
# serial number.
$proximity.Range("A53:A63").NumberFormat = "@"

$proximityRows = @(
    @("2026-02-01","15:16:19","15:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-02-01","15:16:27","15:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-02-01","15:16:34","15:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-02-01","15:16:39","15:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-02-01","15:16:45","15:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-02-01","15:16:57","15:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-02-01","15:16:59","15:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-02-01","15:17:02","15:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-02-01","15:17:04","15:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-02-01","15:17:08","15:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-02-01","15:17:14","15:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door")
)

$r = 53
foreach ($row in $proximityRows) {
    $proximity.Cells.Item($r, 1).Value = $row[0]
    $proximity.Cells.Item($r, 2).Value = $row[1]
    $proximity.Cells.Item($r, 3).Value = $row[2]
    $proximity.Cells.Item($r, 4).Value = $row[3]
    $proximity.Cells.Item($r, 5).Value = $row[4]
    $proximity.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "Camera": append rows 18-25 (new "Image Captured" log entries)
# ---------------------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$camera.Range("A18:A25").NumberFormat = "@"

$cameraRows = @(
    @("2026-02-01","15:16:19","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:16:28","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:16:34","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:16:40","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:16:45","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:16:59","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:17:04","15:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","15:17:14","15:00","Living Room Main Door","Image Captured","Active")
)

$r = 18
foreach ($row in $cameraRows) {
    $camera.Cells.Item($r, 1).Value = $row[0]
    $camera.Cells.Item($r, 2).Value = $row[1]
    $camera.Cells.Item($r, 3).Value = $row[2]
    $camera.Cells.Item($r, 4).Value = $row[3]
    $camera.Cells.Item($r, 5).Value = $row[4]
    $camera.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

Write-Output "Proximity now A1:F$($proximity.Cells.Item(63,6).Row)"
Write-Output "Camera now A1:F$($camera.Cells.Item(25,6).Row)"
